$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H33").Value = 382.75
$ws.Range("I33").Value = 391.3125
$ws.Range("K33").Value = 391.3125
$ws.Range("M33").Value = -162.3125
$ws.Range("H69").Value = 7934.6
$ws.Range("I69").Value = 6750
$ws.Range("K69").Value = 20250
$ws.Range("M69").Value = -19376
$ws.Range("H72").Value = 7934.6
$ws.Range("I72").Value = 6750
$ws.Range("K72").Value = 60750
$ws.Range("M72").Value = -56382
$ws.Range("H106").Value = 20327
$ws.Range("I106").Value = 22675.3
$ws.Range("K106").Value = 22675.3
$ws.Range("M106").Value = -22044.3
$ws.Range("H115").Value = 335.36365
$ws.Range("I115").Value = 346.66666
$ws.Range("J115").Value = 284.5
$ws.Range("K115").Value = 1039.99998
$ws.Range("L115").Value = 853.5
$ws.Range("M115").Value = 527.0000199999999
$ws.Range("N115").Value = -3987.5
$ws.Range("H136").Value = 175492.25
$ws.Range("J136").Value = 175492.25
$ws.Range("L136").Value = 175492.25
$ws.Range("N136").Value = -185692.25
$ws.Range("H138").Value = 4469.08
$ws.Range("I138").Value = 2839
$ws.Range("J138").Value = 4650.2
$ws.Range("K138").Value = 8517
$ws.Range("L138").Value = 13950.6
$ws.Range("M138").Value = -3377
$ws.Range("N138").Value = -24230.6

$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 1436.8534
$ws.Range("I32").Value = 1291.9429
$ws.Range("K32").Value = 1291.9429
$ws.Range("M32").Value = -1004.9429
$ws.Range("H61").Value = 3898.25
$ws.Range("I61").Value = 3798.0908
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 3798.0908
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -3586.0908
$ws.Range("N61").Value = -5424
$ws.Range("H74").Value = 42300.285
$ws.Range("I74").Value = 8935.130999999999
$ws.Range("J74").Value = 195780
$ws.Range("K74").Value = 8935.130999999999
$ws.Range("L74").Value = 195780
$ws.Range("M74").Value = -8061.130999999999
$ws.Range("N74").Value = -197528
$ws.Range("H77").Value = 42300.285
$ws.Range("I77").Value = 8935.130999999999
$ws.Range("J77").Value = 195780
$ws.Range("K77").Value = 44675.655
$ws.Range("L77").Value = 978900
$ws.Range("M77").Value = -40307.655
$ws.Range("N77").Value = -987636
$ws.Range("H80").Value = 50447.5
$ws.Range("J80").Value = 50447.5
$ws.Range("L80").Value = 50447.5
$ws.Range("N80").Value = -52443.5
$ws.Range("H83").Value = 50447.5
$ws.Range("J83").Value = 50447.5
$ws.Range("L83").Value = 151342.5
$ws.Range("N83").Value = -161326.5
$ws.Range("H132").Value = 2522.919
$ws.Range("I132").Value = 2020.1428
$ws.Range("K132").Value = 6060.428400000001
$ws.Range("M132").Value = -3530.428400000001
$ws.Range("H136").Value = 3898.25
$ws.Range("I136").Value = 3798.0908
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 11394.2724
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -8844.2724
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item(3)
$ws.Range("H63").Value = 37447.5
$ws.Range("J63").Value = 37447.5
$ws.Range("L63").Value = 37447.5
$ws.Range("N63").Value = -38819.5
$ws.Range("H66").Value = 37447.5
$ws.Range("J66").Value = 37447.5
$ws.Range("L66").Value = 112342.5
$ws.Range("N66").Value = -119206.5
$ws.Range("H105").Value = 54191.367
$ws.Range("I105").Value = 84761.75
$ws.Range("J105").Value = 1785
$ws.Range("K105").Value = 84761.75
$ws.Range("L105").Value = 1785
$ws.Range("M105").Value = -83014.75
$ws.Range("N105").Value = -5279

$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 19834.691
$ws.Range("J31").Value = 55952.39
$ws.Range("L31").Value = 55952.39
$ws.Range("N31").Value = -56542.39
$ws.Range("H34").Value = 19834.691
$ws.Range("J34").Value = 55952.39
$ws.Range("L34").Value = 55952.39
$ws.Range("N34").Value = -56356.39
$ws.Range("H58").Value = 2873.5715
$ws.Range("I58").Value = 2346.0667
$ws.Range("K58").Value = 2346.0667
$ws.Range("M58").Value = -2143.0667
$ws.Range("H62").Value = 2244.4
$ws.Range("I62").Value = 1740.6666
$ws.Range("K62").Value = 1740.6666
$ws.Range("M62").Value = -1116.6666
$ws.Range("H65").Value = 2244.4
$ws.Range("I65").Value = 1740.6666
$ws.Range("K65").Value = 8703.333000000001
$ws.Range("M65").Value = -5583.333000000001
$ws.Range("H94").Value = 4835.4
$ws.Range("I94").Value = 3150
$ws.Range("K94").Value = 3150
$ws.Range("M94").Value = -2699
$ws.Range("H122").Value = 2193.36
$ws.Range("I122").Value = 1533.8125
$ws.Range("K122").Value = 4601.4375
$ws.Range("M122").Value = -2151.4375
$ws.Range("H134").Value = 3454.6667
$ws.Range("J134").Value = 4799
$ws.Range("L134").Value = 14397
$ws.Range("N134").Value = -19467
$ws.Range("H136").Value = 2873.5715
$ws.Range("I136").Value = 2346.0667
$ws.Range("K136").Value = 7038.2001
$ws.Range("M136").Value = -4488.2001

$ws = $wb.Worksheets.Item(5)
$ws.Range("H56").Value = 6849.933
$ws.Range("I56").Value = 6849.933
$ws.Range("K56").Value = 6849.933
$ws.Range("M56").Value = -6319.933
$ws.Range("H113").Value = 2275.6667
$ws.Range("J113").Value = 1623.75
$ws.Range("L113").Value = 4871.25
$ws.Range("N113").Value = -9211.25
$ws.Range("H131").Value = 2654.111
$ws.Range("J131").Value = 3446.1667
$ws.Range("L131").Value = 10338.5001
$ws.Range("N131").Value = -20418.5001

$ws = $wb.Worksheets.Item(6)
$ws.Range("H14").Value = 14290400
$ws.Range("I14").Value = 3619.8
$ws.Range("J14").Value = 50007350
$ws.Range("K14").Value = 3619.8
$ws.Range("L14").Value = 50007350
$ws.Range("M14").Value = -3451.8
$ws.Range("N14").Value = -50007686
$ws.Range("H63").Value = 39995
$ws.Range("J63").Value = 39995
$ws.Range("L63").Value = 39995
$ws.Range("N63").Value = -41367
$ws.Range("H66").Value = 39995
$ws.Range("J66").Value = 39995
$ws.Range("L66").Value = 119985
$ws.Range("N66").Value = -126849
$ws.Range("H70").Value = 7049.375
$ws.Range("J70").Value = 7450
$ws.Range("L70").Value = 7450
$ws.Range("N70").Value = -7990
$ws.Range("H73").Value = 7049.375
$ws.Range("J73").Value = 7450
$ws.Range("L73").Value = 7450
$ws.Range("N73").Value = -9322
$ws.Range("H86").Value = 41663.332
$ws.Range("J86").Value = 41663.332
$ws.Range("L86").Value = 41663.332
$ws.Range("N86").Value = -44035.332
$ws.Range("H89").Value = 41663.332
$ws.Range("J89").Value = 41663.332
$ws.Range("L89").Value = 124989.996
$ws.Range("N89").Value = -136845.996
$ws.Range("H102").Value = 42822.08
$ws.Range("I102").Value = 1706.1333
$ws.Range("K102").Value = 1706.1333
$ws.Range("M102").Value = -84.13329999999996
$ws.Range("H132").Value = 3094.5715
$ws.Range("I132").Value = 3004.6843
$ws.Range("K132").Value = 9014.052899999999
$ws.Range("M132").Value = -6484.052899999999

$ws = $wb.Worksheets.Item(7)
$ws.Range("H2").Value = 10004950
$ws.Range("I2").Value = 20000334
$ws.Range("J2").Value = 9566.666999999999
$ws.Range("K2").Value = 20000334
$ws.Range("L2").Value = 9566.666999999999
$ws.Range("M2").Value = -20000222
$ws.Range("N2").Value = -9790.666999999999
$ws.Range("H46").Value = 5920.552
$ws.Range("J46").Value = 6681.5415
$ws.Range("L46").Value = 6681.5415
$ws.Range("N46").Value = -7057.5415
$ws.Range("H50").Value = 37580
$ws.Range("I50").Value = 30076
$ws.Range("J50").Value = 45084
$ws.Range("K50").Value = 30076
$ws.Range("L50").Value = 45084
$ws.Range("M50").Value = -29439
$ws.Range("N50").Value = -46358
$ws.Range("H122").Value = 6751.4
$ws.Range("J122").Value = 7681.077
$ws.Range("L122").Value = 23043.231
$ws.Range("N122").Value = -27943.231
$ws.Range("H132").Value = 5571.6206
$ws.Range("I132").Value = 4919.625
$ws.Range("K132").Value = 14758.875
$ws.Range("M132").Value = -12228.875

$ws = $wb.Worksheets.Item(8)
$ws.Range("H13").Value = 3306.077
$ws.Range("I13").Value = 2997.4
$ws.Range("J13").Value = 4335
$ws.Range("K13").Value = 2997.4
$ws.Range("L13").Value = 4335
$ws.Range("M13").Value = -2857.4
$ws.Range("N13").Value = -4615
$ws.Range("H81").Value = 2400
$ws.Range("I81").Value = 2600
$ws.Range("K81").Value = 5200
$ws.Range("M81").Value = -4139
$ws.Range("H84").Value = 2400
$ws.Range("I84").Value = 2600
$ws.Range("K84").Value = 26000
$ws.Range("M84").Value = -20696
$ws.Range("H126").Value = 4075.7646
$ws.Range("I126").Value = 3437.0908
$ws.Range("J126").Value = 5246.6665
$ws.Range("K126").Value = 10311.2724
$ws.Range("L126").Value = 15739.9995
$ws.Range("M126").Value = -7841.2724
